$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All target cells are text (mirrors the source workbook, where B/C/D/E
# are stored as inlineStr). Force text format per-cell before writing so
# numeric-looking strings (e.g. "68.368.58", "0.0000291") are not coerced
# into numbers by Excels automatic type inference.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D2").Value = "68.368.58"
$ws.Range("E2").Value = "  +0.73%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("D3").Value = "3.622.77"
$ws.Range("E3").Value = "  -0.81%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.48%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("D5").Value = "585.28"
$ws.Range("E5").Value = "  -0.72%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("D6").Value = "194.80"
$ws.Range("E6").Value = "  +3.47%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("D7").Value = "3.616.71"
$ws.Range("E7").Value = "  -0.75%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("D8").Value = "0.622"
$ws.Range("E8").Value = "  +1.07%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("D10").Value = "0.681"
$ws.Range("E10").Value = "  -0.95%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.32%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("D12").Value = "55.86"
$ws.Range("E12").Value = "  +0.61%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000291"
$ws.Range("E13").Value = "  +11.03%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("D14").Value = "10.09"
$ws.Range("E14").Value = "  -0.15%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D15").Value = "4.190.85"
$ws.Range("E15").Value = "  -1.03%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("D16").Value = "3.616.36"
$ws.Range("E16").Value = "  -1.06%  "

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.10%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("D18").Value = "12.56"
$ws.Range("E18").Value = "  +1.29%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("D19").Value = "68.199.78"
$ws.Range("E19").Value = "  +0.66%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("D20").Value = "18.60"
$ws.Range("E20").Value = "  -0.21%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.91%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("D22").Value = "405.24"
$ws.Range("E22").Value = "  +1.66%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("D23").Value = "13.36"
$ws.Range("E23").Value = "  +24.87%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("D24").Value = "4.27"
$ws.Range("E24").Value = "  -2.18%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("D25").Value = "86.27"
$ws.Range("E25").Value = "  -0.58%  "

$ws.Range("B26").NumberFormat = "@"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "4.02"
$ws.Range("E26").Value = "  +10.20%  "

$ws.Range("B27").NumberFormat = "@"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D27").Value = "2.96"
$ws.Range("E27").Value = "  +2.04%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("D28").Value = "12.66"
$ws.Range("E28").Value = "  +2.32%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("D29").Value = "6.16"
$ws.Range("E29").Value = "  +1.32%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("D30").Value = "8.16"
$ws.Range("E30").Value = "  +15.80%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("D31").Value = "9.20"
$ws.Range("E31").Value = "  +0.39%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("D32").Value = "31.74"
$ws.Range("E32").Value = "  +0.58%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("D33").Value = "693.30"
$ws.Range("E33").Value = "  +14.48%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("D34").Value = "12.30"
$ws.Range("E34").Value = "  +1.72%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.71%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("D36").Value = "64.85"
$ws.Range("E36").Value = "  -3.92%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("D37").Value = "42.87"
$ws.Range("E37").Value = "  -1.04%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("D38").Value = "0.420"
$ws.Range("E38").Value = "  +8.96%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.04%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0795"
$ws.Range("E40").Value = "  +7.76%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("D41").Value = "2.91"
$ws.Range("E41").Value = "  +18.01%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D42").Value = "3.15"
$ws.Range("E42").Value = "  +10.84%  "

$ws.Range("B43").NumberFormat = "@"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "0.135"
$ws.Range("E43").Value = "  -0.35%  "

$ws.Range("B44").NumberFormat = "@"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "3.162.59"
$ws.Range("E44").Value = "  +15.58%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.42%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0426"
$ws.Range("E46").Value = "  +2.14%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.13%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("D48").Value = "8.89"
$ws.Range("E48").Value = "  +1.57%  "

$ws.Range("B49").NumberFormat = "@"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "143.66"
$ws.Range("E49").Value = "  +1.14%  "

$ws.Range("B50").NumberFormat = "@"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("B50").Value = "ApeXProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D50").Value = "3.12"
$ws.Range("E50").Value = "  -1.91%  "

$ws.Range("B51").NumberFormat = "@"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("B51").Value = "WEMIXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D51").Value = "2.61"
$ws.Range("E51").Value = "  +1.07%  "
